$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clone row 8's formatting down onto the 4 new rows (9-12) first, so the
#     number formats (in particular the "Text" format used by columns D/E)
#     are already in place before any values are typed in - this keeps
#     leading zeros on phone numbers intact. ---
$ws.Range("A8:I8").Copy()
$ws.Range("A9:I12").PasteSpecial(-4122)   # xlPasteFormats

# Rows 10-11 (D:E) use the "wrap text" variant that already exists on row 3
$ws.Range("D3:E3").Copy()
$ws.Range("D10:E11").PasteSpecial(-4122)  # xlPasteFormats

# Match the row height used by the rest of the table
$ws.Range("9:12").RowHeight = 15

# --- Row 9 : SITHA RAMADHANI A. ---
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "SITHA RAMADHANI A."
$ws.Cells.Item(9, 3).Value = 1010202208
$ws.Cells.Item(9, 4).Value = "087623918732"
$ws.Cells.Item(9, 5).Value = "0891278318723"
$ws.Cells.Item(9, 6).Value = "ramasd123@gmail.com"
$ws.Cells.Item(9, 7).Value = "KOTA CIMAHI"
$ws.Cells.Item(9, 8).Value = "KEPERAWATAN"
$ws.Cells.Item(9, 9).Value = "STIKES JIWA SEHAT"

# --- Row 10 : ARIF HAKIM ---
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "ARIF HAKIM"
$ws.Cells.Item(10, 3).Value = 1010202209
$ws.Cells.Item(10, 4).Value = "0822891238911"
$ws.Cells.Item(10, 5).Value = "0822891238911"
$ws.Cells.Item(10, 6).Value = "arif_hakim@gmail.com"
$ws.Cells.Item(10, 7).Value = "KOTA CIMAHI"
$ws.Cells.Item(10, 8).Value = "KEPERAWATAN"
$ws.Cells.Item(10, 9).Value = "STIKES JIWA SEHAT"

# --- Row 11 : ADI HARDIANSYAH ---
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "ADI HARDIANSYAH"
$ws.Cells.Item(11, 3).Value = 1010202210
$ws.Cells.Item(11, 4).Value = "0889263223683"
$ws.Cells.Item(11, 5).Value = "0889263223683"
$ws.Cells.Item(11, 6).Value = "adihhardianr@gmail.com"
$ws.Cells.Item(11, 7).Value = "KAB. BANDUNG BARAT"
$ws.Cells.Item(11, 8).Value = "KEPERAWATAN"
$ws.Cells.Item(11, 9).Value = "STIKES JIWA SEHAT"

# --- Row 12 : NANDANG ---
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "NANDANG"
$ws.Cells.Item(12, 3).Value = 1010202211
$ws.Cells.Item(12, 4).Value = "0898298332323"
$ws.Cells.Item(12, 5).Value = "0898298332323"
$ws.Cells.Item(12, 6).Value = "nandang@gmail.com"
$ws.Cells.Item(12, 7).Value = "KAB. BANDUNG BARAT"
$ws.Cells.Item(12, 8).Value = "KEPERAWATAN"
$ws.Cells.Item(12, 9).Value = "STIKES JIWA SEHAT"

# --- Hyperlinks for the new e-mail cells ---
$ws.Hyperlinks.Add($ws.Range("F9"), "mailto:ramasd123@gmail.com")
$ws.Hyperlinks.Add($ws.Range("F10"), "mailto:arif_hakim@gmail.com")
$ws.Hyperlinks.Add($ws.Range("F11"), "mailto:adihhardianr@gmail.com")
$ws.Hyperlinks.Add($ws.Range("F12"), "mailto:nandang@gmail.com")

# Hyperlinks.Add stamps column F with the built-in blue/underlined
# "Hyperlink" style; restore the sheet's own (non-underlined) look that
# the rest of the e-mail column already uses.
$ws.Range("F8").Copy()
$ws.Range("F9:F12").PasteSpecial(-4122)   # xlPasteFormats

# --- Misc view state ---
$ws.Range("D1").Select()
